$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1257.6389
$ws.Range("J129").Value = 1278.6765
$ws.Range("L129").Value = 3836.0295
$ws.Range("N129").Value = -13836.0295
$ws.Range("H134").Value = 116161.9
$ws.Range("J134").Value = 116161.9
$ws.Range("L134").Value = 116161.9
$ws.Range("N134").Value = -126301.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4004.3333
$ws.Range("J61").Value = 4517.75
$ws.Range("L61").Value = 4517.75
$ws.Range("N61").Value = -4941.75
$ws.Range("H74").Value = 1423.75
$ws.Range("I74").Value = 1385.2
$ws.Range("J74").Value = 1616.5
$ws.Range("K74").Value = 1385.2
$ws.Range("L74").Value = 1616.5
$ws.Range("M74").Value = -511.2
$ws.Range("N74").Value = -3364.5
$ws.Range("H77").Value = 1423.75
$ws.Range("I77").Value = 1385.2
$ws.Range("J77").Value = 1616.5
$ws.Range("K77").Value = 6926
$ws.Range("L77").Value = 8082.5
$ws.Range("M77").Value = -2558
$ws.Range("N77").Value = -16818.5
$ws.Range("H132").Value = 2445.3333
$ws.Range("I132").Value = 2150.7932
$ws.Range("J132").Value = 3299.5
$ws.Range("K132").Value = 6452.3796
$ws.Range("L132").Value = 9898.5
$ws.Range("M132").Value = -3922.3796
$ws.Range("N132").Value = -14958.5
$ws.Range("H134").Value = 55255.668
$ws.Range("J134").Value = 55255.668
$ws.Range("L134").Value = 55255.668
$ws.Range("N134").Value = -65395.668
$ws.Range("H136").Value = 4004.3333
$ws.Range("J136").Value = 4517.75
$ws.Range("L136").Value = 13553.25
$ws.Range("N136").Value = -18653.25
$ws.Range("H139").Value = 96632
$ws.Range("J139").Value = 96632
$ws.Range("L139").Value = 96632
$ws.Range("N139").Value = -106912
$ws.Range("H141").Value = 48122.715
$ws.Range("J141").Value = 48122.715
$ws.Range("L141").Value = 48122.715
$ws.Range("N141").Value = -58482.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1571.825
$ws.Range("I132").Value = 1521.5278
$ws.Range("J132").Value = 2024.5
$ws.Range("K132").Value = 4564.5834
$ws.Range("L132").Value = 6073.5
$ws.Range("M132").Value = -2034.5834
$ws.Range("N132").Value = -11133.5
$ws.Range("H140").Value = 70262
$ws.Range("J140").Value = 77827.5
$ws.Range("L140").Value = 77827.5
$ws.Range("N140").Value = -88187.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7951.8125
$ws.Range("I122").Value = 593
$ws.Range("K122").Value = 5337
$ws.Range("M122").Value = -2887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4650.2
$ws.Range("I41").Value = 1417
$ws.Range("J41").Value = 9500
$ws.Range("K41").Value = 1417
$ws.Range("L41").Value = 9500
$ws.Range("M41").Value = -1062
$ws.Range("N41").Value = -10210
$ws.Range("H132").Value = 1771.9395
$ws.Range("I132").Value = 1502.8636
$ws.Range("J132").Value = 2310.0908
$ws.Range("K132").Value = 4508.5908
$ws.Range("L132").Value = 6930.2724
$ws.Range("M132").Value = -1978.5908
$ws.Range("N132").Value = -11990.2724
$ws.Range("H138").Value = 45967.9
$ws.Range("J138").Value = 45967.9
$ws.Range("L138").Value = 45967.9
$ws.Range("N138").Value = -56247.9
$ws.Range("H141").Value = 60507.5
$ws.Range("J141").Value = 60507.5
$ws.Range("L141").Value = 60507.5
$ws.Range("N141").Value = -70867.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8378.944
$ws.Range("I132").Value = 8401.235000000001
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 25203.705
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -22673.705
$ws.Range("N132").Value = -29060
$ws.Range("H133").Value = 65354.383
$ws.Range("J133").Value = 65354.383
$ws.Range("L133").Value = 65354.383
$ws.Range("N133").Value = -70414.383
$ws.Range("H134").Value = 38689
$ws.Range("J134").Value = 38689
$ws.Range("L134").Value = 38689
$ws.Range("N134").Value = -48829
$ws.Range("H136").Value = 6345.6924
$ws.Range("I136").Value = 6832.524
$ws.Range("J136").Value = 4301
$ws.Range("K136").Value = 20497.572
$ws.Range("L136").Value = 12903
$ws.Range("M136").Value = -17947.572
$ws.Range("N136").Value = -18003
$ws.Range("H137").Value = 78416
$ws.Range("J137").Value = 78416
$ws.Range("L137").Value = 78416
$ws.Range("N137").Value = -88616
$ws.Range("H140").Value = 84150
$ws.Range("J140").Value = 84150
$ws.Range("L140").Value = 84150
$ws.Range("N140").Value = -94510
$ws.Range("H141").Value = 55750
$ws.Range("J141").Value = 55750
$ws.Range("L141").Value = 55750
$ws.Range("N141").Value = -66110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 52736.168
$ws.Range("J46").Value = 52736.168
$ws.Range("L46").Value = 52736.168
$ws.Range("N46").Value = -53198.168
$ws.Range("H132").Value = 1883.909
$ws.Range("I132").Value = 2073
$ws.Range("J132").Value = 1593
$ws.Range("K132").Value = 6219
$ws.Range("L132").Value = 4779
$ws.Range("M132").Value = -3689
$ws.Range("N132").Value = -9839
$ws.Range("H134").Value = 52736.168
$ws.Range("J134").Value = 52736.168
$ws.Range("L134").Value = 158208.504
$ws.Range("N134").Value = -163278.504
$ws.Range("H135").Value = 70271.664
$ws.Range("J135").Value = 70271.664
$ws.Range("L135").Value = 70271.664
$ws.Range("N135").Value = -80411.664
$ws.Range("H136").Value = 2083.1177
$ws.Range("I136").Value = 2087.7307
$ws.Range("J136").Value = 2068.125
$ws.Range("K136").Value = 6263.1921
$ws.Range("L136").Value = 6204.375
$ws.Range("M136").Value = -3713.1921
$ws.Range("N136").Value = -11304.375
$ws.Range("H137").Value = 67017.27
$ws.Range("J137").Value = 67017.27
$ws.Range("L137").Value = 67017.27
$ws.Range("N137").Value = -77217.27
$ws.Range("H139").Value = 62686.54
$ws.Range("J139").Value = 62686.54
$ws.Range("L139").Value = 62686.54
$ws.Range("N139").Value = -72966.54000000001
$ws.Range("H140").Value = 38216.332
$ws.Range("J140").Value = 38216.332
$ws.Range("L140").Value = 38216.332
$ws.Range("N140").Value = -48576.332
$ws.Range("H141").Value = 63344.09
$ws.Range("J141").Value = 63344.09
$ws.Range("L141").Value = 63344.09
$ws.Range("N141").Value = -73704.09
